$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 <-> Row 30: swap Id (A), Aktivitet (M), Ost (Q), Nord (R)
$ws.Range("A28").Value = 111985427
$ws.Range("M28").Value = "äldre spår"
$ws.Range("Q28").Value = 562013.8588788129
$ws.Range("R28").Value = 7307395.432326685

$ws.Range("A30").Value = 111985426
$ws.Range("M30").Value = "färska spår"
$ws.Range("Q30").Value = 562016.495064693
$ws.Range("R30").Value = 7307406.130124222

# Row 29 <-> Row 31: swap Id (A), Taxonsorteringsordning (B), Rödlistade (D),
# TaxonId (E), Artnamn (F), Vetenskapligt namn (G), Auktor (H), Ost (Q), Nord (R)
$ws.Range("A29").Value = 111959833
$ws.Range("B29").Value = 81248
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 1312
$ws.Range("F29").Value = "Gammelgransskål"
$ws.Range("G29").Value = "Pseudographis pinicola"
$ws.Range("H29").Value = "(Nyl.) Rehm"
$ws.Range("Q29").Value = 562008.043958648
$ws.Range("R29").Value = 7307380.163763028

$ws.Range("A31").Value = 111959825
$ws.Range("B31").Value = 89747
$ws.Range("D31").Value = "VU"
$ws.Range("E31").Value = 2063
$ws.Range("F31").Value = "Grantickeporing"
$ws.Range("G31").Value = "Skeletocutis chrysella"
$ws.Range("H31").Value = "Niemelä"
$ws.Range("Q31").Value = 561809.3860941484
$ws.Range("R31").Value = 7307206.837266683

# Rows 41 -> 42 -> 43 -> 41 cyclic rotation of:
# Id (A), Taxonsorteringsordning (B), TaxonId (E), Artnamn (F),
# Vetenskapligt namn (G), Auktor (H), Ost (Q), Nord (R)
$ws.Range("A41").Value = 111959828
$ws.Range("B41").Value = 89423
$ws.Range("E41").Value = 5432
$ws.Range("F41").Value = "Granticka"
$ws.Range("G41").Value = "Porodaedalea chrysoloma"
$ws.Range("H41").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q41").Value = 561996.2119675929
$ws.Range("R41").Value = 7307377.861109382

$ws.Range("A42").Value = 111959829
$ws.Range("B42").Value = 77515
$ws.Range("E42").Value = 6425
$ws.Range("F42").Value = "Garnlav"
$ws.Range("G42").Value = "Alectoria sarmentosa"
$ws.Range("H42").Value = "(Ach.) Ach."

$ws.Range("A43").Value = 111959818
$ws.Range("B43").Value = 89686
$ws.Range("E43").Value = 658
$ws.Range("F43").Value = "Rosenticka"
$ws.Range("G43").Value = "Rhodofomes roseus"
$ws.Range("H43").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q43").Value = 561893.9245207607
$ws.Range("R43").Value = 7307219.714951258
